$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update betting odds values for 2025-05-15 games (rows 12, 13, 15, 17, 19, 23)

# Row 12
$ws.Range("G12").Value = 1.85
$ws.Range("I12").Value = 4.5
$ws.Range("M12").Value = 2.18
$ws.Range("P12").Value = 1.6
$ws.Range("Q12").Value = 2.07
$ws.Range("T12").Value = 4.8
$ws.Range("U12").Value = 7.1
$ws.Range("Z12").Value = 5.9
$ws.Range("AA12").Value = 6.3
$ws.Range("AF12").Value = 23
$ws.Range("AG12").Value = 16.5

# Row 13
$ws.Range("J13").Value = 1.07
$ws.Range("K13").Value = 9

# Row 15
$ws.Range("G15").Value = 4.9
$ws.Range("H15").Value = 3.45
$ws.Range("I15").Value = 1.62
$ws.Range("N15").Value = 1.98
$ws.Range("O15").Value = 1.65
$ws.Range("T15").Value = 10
$ws.Range("U15").Value = 22
$ws.Range("V15").Value = 13
$ws.Range("W15").Value = 65
$ws.Range("X15").Value = 40
$ws.Range("Z15").Value = 8.5
$ws.Range("AA15").Value = 6
$ws.Range("AC15").Value = 70
$ws.Range("AE15").Value = 5.1
$ws.Range("AF15").Value = 6
$ws.Range("AG15").Value = 7
$ws.Range("AH15").Value = 9.75
$ws.Range("AI15").Value = 11.5
$ws.Range("AJ15").Value = 24

# Row 17
$ws.Range("H17").Value = 3.2
$ws.Range("I17").Value = 3.5
$ws.Range("J17").Value = 1.07
$ws.Range("K17").Value = 9
$ws.Range("L17").Value = 1.33
$ws.Range("M17").Value = 3.25
$ws.Range("N17").Value = 2.08
$ws.Range("O17").Value = 1.73
$ws.Range("P17").Value = 1.44
$ws.Range("Q17").Value = 2.63
$ws.Range("R17").Value = 1.8
$ws.Range("S17").Value = 1.95
$ws.Range("T17").Value = 7.5
$ws.Range("U17").Value = 10
$ws.Range("V17").Value = 9.5
$ws.Range("X17").Value = 19
$ws.Range("Y17").Value = 29
$ws.Range("Z17").Value = 8.5
$ws.Range("AA17").Value = 6
$ws.Range("AB17").Value = 15
$ws.Range("AC17").Value = 51
$ws.Range("AD17").Value = 251
$ws.Range("AE17").Value = 9.5
$ws.Range("AI17").Value = 29
$ws.Range("AJ17").Value = 41

# Row 19
$ws.Range("K19").Value = 13

# Row 23
$ws.Range("G23").Value = 3.45
$ws.Range("I23").Value = 2.07
$ws.Range("N23").Value = 1.95
$ws.Range("R23").Value = 1.75
$ws.Range("S23").Value = 1.95
$ws.Range("T23").Value = 10
$ws.Range("V23").Value = 12.5
$ws.Range("W23").Value = 50
$ws.Range("X23").Value = 32
$ws.Range("Y23").Value = 40
$ws.Range("AE23").Value = 7.2
$ws.Range("AG23").Value = 9.25
$ws.Range("AI23").Value = 18
